# Insert a new data row at row 110 (pushes the existing rows 110..239 down
# to 111..240) and populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 44539
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100112023
$ws.Range("G110").Value = "Brócoli"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 550
$ws.Range("L110").Value = 550
$ws.Range("M110").Value = 550
$ws.Range("N110").Value = "$/unidad"
$ws.Range("O110").Value = "Región del Maule"
$ws.Range("P110").Value = 550
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"
